$d = $word.ActiveDocument

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph, which is
# immediately preceded by a blank paragraph and immediately followed by the
# site footer/copyright paragraph. The commit removes all three paragraphs,
# leaving the previous content paragraph directly followed by the paragraph
# that used to come after the footer.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 1) {
    $startPara = $d.Paragraphs.Item($targetIndex - 1)
    $endPara = $d.Paragraphs.Item($targetIndex + 1)
    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}
